$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsScripts = $wb.Worksheets.Item("TestScripts")
$wsLoginData = $wb.Worksheets.Item("LoginData")
$wsMaddies = $wb.Worksheets.Item("Maddieslogindata")

# --- TestScripts: flip pass/fail statuses ---
$wsScripts.Range("B2").Value = "FAIL"
$wsScripts.Range("B3").Value = "Pass"
$wsScripts.Range("B4").Value = "FAIL"

# --- LoginData: same test email for every row, passwords stored as text ---
$wsLoginData.Range("B2").Value = "1529830551216@gmail.com"
$wsLoginData.Range("B3").Value = "1529830551216@gmail.com"
$wsLoginData.Range("B4").Value = "1529830551216@gmail.com"

# Match the "Text" number format cell style already used on the Maddieslogindata sheet
# (same border/font, numFmtId 49) before writing the literal "password" text.
$wsMaddies.Range("C2").Copy()
$wsLoginData.Range("C2:C4").PasteSpecial(-4122)

$wsLoginData.Range("C2").Value = "password"
$wsLoginData.Range("C3").Value = "password"
$wsLoginData.Range("C4").Value = "password"

# Remove the mailto hyperlinks from LoginData (still present on Maddieslogindata)
$wsLoginData.Range("B2").Hyperlinks.Delete()

# --- Selections / active sheet ---
$wsScripts.Activate()
$wsScripts.Range("A4").Select()

$wsLoginData.Activate()
$wsLoginData.Range("B2:B4").Select()

$wsMaddies.Activate()
$wsMaddies.Range("C2:C4").Select()

$wsScripts.Activate()
